$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 125; this shifts all existing rows
# 125..179 down to 126..180 (Excel copies formatting, e.g. the date
# style on column D, from the row above automatically).
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row 125 with the new weekly record.
$ws.Cells.Item(125, 1).Value = 8
$ws.Cells.Item(125, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(125, 3).Value = "Coquimbo"
$ws.Cells.Item(125, 4).Value = 44875
$ws.Cells.Item(125, 5).Value = 4
$ws.Cells.Item(125, 6).Value = 100112044
$ws.Cells.Item(125, 7).Value = "Perejil"
$ws.Cells.Item(125, 8).Value = "Sin especificar"
$ws.Cells.Item(125, 9).Value = "Primera"
$ws.Cells.Item(125, 10).Value = 2560
$ws.Cells.Item(125, 11).Value = 1500
$ws.Cells.Item(125, 12).Value = 2000
$ws.Cells.Item(125, 13).Value = 1750
$ws.Cells.Item(125, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(125, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(125, 16).Value = 1167
$ws.Cells.Item(125, 17).Value = 1.5
$ws.Cells.Item(125, 18).Value = "Hortaliza"
